$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 644; existing rows 644-676 shift down to 645-677
$ws.Rows.Item(644).Insert()

# Populate the new row 644 with the new entry's data
$ws.Cells.Item(644, 1).Value = 3
$ws.Cells.Item(644, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(644, 3).Value = "Coquimbo"
$ws.Cells.Item(644, 4).Value = 45267
$ws.Cells.Item(644, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(644, 5).Value = 5
$ws.Cells.Item(644, 6).Value = 100112040
$ws.Cells.Item(644, 7).Value = "Cilantro"
$ws.Cells.Item(644, 8).Value = "Sin especificar"
$ws.Cells.Item(644, 9).Value = "Primera"
$ws.Cells.Item(644, 10).Value = 180
$ws.Cells.Item(644, 11).Value = 7500
$ws.Cells.Item(644, 12).Value = 8000
$ws.Cells.Item(644, 13).Value = 7778
$ws.Cells.Item(644, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(644, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(644, 16).Value = 2593
$ws.Cells.Item(644, 17).Value = 3
$ws.Cells.Item(644, 18).Value = "Hortaliza"
